# Update "想去人数" (want-to-go count) figures in both the "展览" and
# "全部类型" worksheets, matching the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - row number => new value for column F
$sheet1Updates = @{
    2  = 6825
    4  = 436
    8  = 109
    9  = 103
    12 = 26
    13 = 179
    14 = 423
    16 = 1627
    17 = 27
    18 = 3437
    19 = 22
    21 = 10
    22 = 2075
    23 = 169
    24 = 5
    28 = 11
}

# Sheet "全部类型" - row number => new value for column F
$sheet4Updates = @{
    2  = 6825
    4  = 436
    9  = 109
    10 = 103
    13 = 26
    14 = 179
    15 = 423
    17 = 1627
    18 = 27
    19 = 3437
    20 = 22
    22 = 10
    23 = 2075
    24 = 169
    25 = 5
    29 = 11
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
